# Fixing cancel count combinor
# Updates raw input values in the summary sheet (formula cells recalc automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Adrian Hardin
$ws.Range("C3").Value = 975

# Row 5 - Cory Caldwell
$ws.Range("M5").Value = 53

# Row 6 - Darrin Neal
$ws.Range("C6").Value = 150
$ws.Range("J6").Value = 0

# Row 7 - David Strehlow
$ws.Range("D7").Value = 4
$ws.Range("M7").Value = 237

# Row 8 - Garrett McKinzie
$ws.Range("C8").Value = 933
$ws.Range("M8").Value = 125

# Row 9 - Howard Seigle
$ws.Range("C9").Value = 388
$ws.Range("M9").Value = 86

# Row 11 - Justin Cohen
$ws.Range("C11").Value = 452
$ws.Range("J11").Value = 0
$ws.Range("M11").Value = 102

# Row 15 - Mark Files
$ws.Range("C15").Value = 736
$ws.Range("D15").Value = 2
$ws.Range("M15").Value = 133

# Row 16 - Matthew Cano
$ws.Range("D16").Value = 10
$ws.Range("M16").Value = 114

# Row 18 - Nick Thompson
$ws.Range("D18").Value = 2
$ws.Range("M18").Value = 140

# Row 19 - Nick Oberle
$ws.Range("M19").Value = 107

# Row 20 - Patrick Stang
$ws.Range("M20").Value = 81

# Row 21 - Pete Scalzo
$ws.Range("C21").Value = 1142
$ws.Range("D21").Value = 6
$ws.Range("M21").Value = 109

# Row 22 - Rhee Cano
$ws.Range("D22").Value = 4
$ws.Range("M22").Value = 61

# Row 23 - Rhi Neal
$ws.Range("C23").Value = 429
$ws.Range("J23").Value = 124
$ws.Range("M23").Value = 114

# Row 24 - Stephanie Kiely
$ws.Range("C24").Value = 388
$ws.Range("J24").Value = 0
$ws.Range("M24").Value = 83

$wb.Save()
